# Trade #70 closed at 2026-02-16 21:35:31 - momentum DOWN +0.000%
#
# This script applies the updates that occurred in a single live-trading
# logging pass:
#   1. Trade #50 on the "leadlag" sheet (row 40) closes out: Exit Price,
#      Status, P&L%, P&L$, Exit Reason and Duration get filled in, and the
#      same closed trade is appended as a new row (51) on "All Trades".
#   2. A brand new trade #70 is opened on the "momentum" sheet and appended
#      as new row 18 there (status OPEN, no exit data yet).
#   3. Summary / Comparison roll-up statistics are refreshed to reflect the
#      now 50 closed leadlag trades.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a text value to a cell while forcing Excel to store it
# literally as text (no autoconversion to number/date/percentage). This
# mirrors typing a leading apostrophe in the Excel UI. The style is reset
# back to Normal afterwards so the quote-prefix marker doesn't linger as
# a visible formatting change on the cell.
# ---------------------------------------------------------------------
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Summary sheet - refresh OVERALL and leadlag roll-up rows
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 50
Set-TextValue $summary.Range("D2") "66.0%"
Set-TextValue $summary.Range("E2") "+12.5527%"
Set-TextValue $summary.Range("F2") "+0.2511%"

Set-TextValue $summary.Range("D3") "45.3%"
Set-TextValue $summary.Range("E3") "+8.3511%"
Set-TextValue $summary.Range("F3") "+0.1576%"

# ---------------------------------------------------------------------
# 2) leadlag sheet - close out trade #50 (row 40)
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Range("G40").Value = 68992.570962
Set-TextValue $leadlag.Range("H40") "CLOSED"
$leadlag.Range("I40").Value = 0.2536
$leadlag.Range("J40").Value = 2.54
Set-TextValue $leadlag.Range("M40") "time_exit_5min"
$leadlag.Range("N40").Value = 5

# ---------------------------------------------------------------------
# 3) momentum sheet - append newly opened trade #70 (row 18)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Range("A18").Value = 70
Set-TextValue $momentum.Range("B18") "2026-02-16"
Set-TextValue $momentum.Range("C18") "21:35:31"
Set-TextValue $momentum.Range("D18") "momentum"
Set-TextValue $momentum.Range("E18") "DOWN"
$momentum.Range("F18").Value = 68518.41499999999
Set-TextValue $momentum.Range("H18") "OPEN"
$momentum.Range("I18").Value = 0
$momentum.Range("J18").Value = 0
$momentum.Range("K18").Value = 0.9
Set-TextValue $momentum.Range("L18") "Downward momentum: -0.291% over 10 samples"
$momentum.Range("N18").Value = 0

# ---------------------------------------------------------------------
# 4) All Trades sheet - append the now-closed trade #50 as row 51
#    (mirrors the leadlag sheet after the close-out above)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("A51").Value = 50
Set-TextValue $allTrades.Range("B51") "2026-02-16"
Set-TextValue $allTrades.Range("C51") "21:30:30"
Set-TextValue $allTrades.Range("D51") "leadlag"
Set-TextValue $allTrades.Range("E51") "UP"
$allTrades.Range("F51").Value = 68818.08
$allTrades.Range("G51").Value = 68992.570962
Set-TextValue $allTrades.Range("H51") "CLOSED"
$allTrades.Range("I51").Value = 0.2536
$allTrades.Range("J51").Value = 2.54
$allTrades.Range("K51").Value = 0.75
Set-TextValue $allTrades.Range("L51") "Binance leading with 0.153% move"
Set-TextValue $allTrades.Range("M51") "time_exit_5min"
$allTrades.Range("N51").Value = 5

# ---------------------------------------------------------------------
# 5) Comparison sheet - refresh leadlag roll-up row
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

Set-TextValue $comparison.Range("C2") "45.3%"
Set-TextValue $comparison.Range("D2") "2.68"
Set-TextValue $comparison.Range("E2") "+0.5552%"
Set-TextValue $comparison.Range("G2") "1.67"
